$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 20
$ws.Range("G2").Value = 16
$ws.Range("I2").Value = $false
